$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 88.85714
$ws.Range("I55").Value = 49.454544
$ws.Range("K55").Value = 49.454544
$ws.Range("M55").Value = 164.545456

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 9323.125
$ws.Range("I62").Value = 13576.667
$ws.Range("J62").Value = 3854.2856
$ws.Range("K62").Value = 13576.667
$ws.Range("L62").Value = 3854.2856
$ws.Range("M62").Value = -12952.667
$ws.Range("N62").Value = -5102.2856

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 9323.125
$ws.Range("I65").Value = 13576.667
$ws.Range("J65").Value = 3854.2856
$ws.Range("K65").Value = 67883.33499999999
$ws.Range("L65").Value = 19271.428
$ws.Range("M65").Value = -64763.33499999999
$ws.Range("N65").Value = -25511.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1387.3704
$ws.Range("I107").Value = 1667.625
$ws.Range("J107").Value = 979.7273
$ws.Range("K107").Value = 1667.625
$ws.Range("L107").Value = 979.7273
$ws.Range("M107").Value = 252.375
$ws.Range("N107").Value = -4819.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1064.6923
$ws.Range("I111").Value = 918.1667
$ws.Range("K111").Value = 2754.5001
$ws.Range("M111").Value = 312.4998999999998

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3552.6538
$ws.Range("J112").Value = 4179.524
$ws.Range("L112").Value = 12538.572
$ws.Range("N112").Value = -14754.572

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 4332.3335
$ws.Range("I113").Value = 3742.5
$ws.Range("J113").Value = 4500.857
$ws.Range("K113").Value = 3742.5
$ws.Range("L113").Value = 4500.857
$ws.Range("M113").Value = -488.5
$ws.Range("N113").Value = -11008.857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4649.75
$ws.Range("I116").Value = 7800
$ws.Range("J116").Value = 4199.7144
$ws.Range("K116").Value = 7800
$ws.Range("L116").Value = 4199.7144
$ws.Range("M116").Value = -4358
$ws.Range("N116").Value = -11083.7144

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1437.5
$ws.Range("I118").Value = 475
$ws.Range("J118").Value = 2400
$ws.Range("K118").Value = 1425
$ws.Range("L118").Value = 7200
$ws.Range("M118").Value = 232
$ws.Range("N118").Value = -10514

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 63411.668
$ws.Range("J139").Value = 70094
$ws.Range("L139").Value = 70094
$ws.Range("N139").Value = -80374

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 407056.1
$ws.Range("I32").Value = 3504.6611
$ws.Range("J32").Value = 2391184
$ws.Range("K32").Value = 3504.6611
$ws.Range("L32").Value = 2391184
$ws.Range("M32").Value = -3217.6611
$ws.Range("N32").Value = -2391758

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2035.7778
$ws.Range("I61").Value = 1797.1578
$ws.Range("J61").Value = 2602.5
$ws.Range("K61").Value = 1797.1578
$ws.Range("L61").Value = 2602.5
$ws.Range("M61").Value = -1585.1578
$ws.Range("N61").Value = -3026.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1620.5927
$ws.Range("I74").Value = 1035.0769
$ws.Range("J74").Value = 2164.2856
$ws.Range("K74").Value = 1035.0769
$ws.Range("L74").Value = 2164.2856
$ws.Range("M74").Value = -161.0769
$ws.Range("N74").Value = -3912.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1620.5927
$ws.Range("I77").Value = 1035.0769
$ws.Range("J77").Value = 2164.2856
$ws.Range("K77").Value = 5175.3845
$ws.Range("L77").Value = 10821.428
$ws.Range("M77").Value = -807.3845000000001
$ws.Range("N77").Value = -19557.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2035.7778
$ws.Range("I136").Value = 1797.1578
$ws.Range("J136").Value = 2602.5
$ws.Range("K136").Value = 5391.4734
$ws.Range("L136").Value = 7807.5
$ws.Range("M136").Value = -2841.4734
$ws.Range("N136").Value = -12907.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2945.3125
$ws.Range("I20").Value = 2187.889
$ws.Range("J20").Value = 3919.1428
$ws.Range("K20").Value = 2187.889
$ws.Range("L20").Value = 3919.1428
$ws.Range("M20").Value = -1940.889
$ws.Range("N20").Value = -4413.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2602.9019
$ws.Range("I86").Value = 2550.2307
$ws.Range("J86").Value = 2657.68
$ws.Range("K86").Value = 2550.2307
$ws.Range("L86").Value = 2657.68
$ws.Range("M86").Value = -1427.2307
$ws.Range("N86").Value = -4903.68

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2602.9019
$ws.Range("I89").Value = 2550.2307
$ws.Range("J89").Value = 2657.68
$ws.Range("K89").Value = 12751.1535
$ws.Range("L89").Value = 13288.4
$ws.Range("M89").Value = -7135.1535
$ws.Range("N89").Value = -24520.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1489.5927
$ws.Range("I99").Value = 1119.7368
$ws.Range("J99").Value = 2368
$ws.Range("K99").Value = 1119.7368
$ws.Range("L99").Value = 2368
$ws.Range("M99").Value = 378.2632000000001
$ws.Range("N99").Value = -5364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 42561.043
$ws.Range("I107").Value = 603.6842
$ws.Range("J107").Value = 201999
$ws.Range("K107").Value = 603.6842
$ws.Range("L107").Value = 201999
$ws.Range("M107").Value = 1316.3158
$ws.Range("N107").Value = -205839

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 28000
$ws.Range("J108").Value = 28000
$ws.Range("L108").Value = 28000
$ws.Range("N108").Value = -35680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4257.357
$ws.Range("I134").Value = 878.2683
$ws.Range("K134").Value = 2634.8049
$ws.Range("M134").Value = -99.80490000000009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 782.73
$ws.Range("I131").Value = 361.46155
$ws.Range("J131").Value = 845.67816
$ws.Range("K131").Value = 1084.38465
$ws.Range("L131").Value = 2537.03448
$ws.Range("M131").Value = 3955.61535
$ws.Range("N131").Value = -12617.03448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 40001256
$ws.Range("I113").Value = 1046.3077
$ws.Range("J113").Value = 83334824
$ws.Range("K113").Value = 1046.3077
$ws.Range("L113").Value = 83334824
$ws.Range("M113").Value = 1123.6923
$ws.Range("N113").Value = -83339164

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2623.5833
$ws.Range("I22").Value = 3254.4443
$ws.Range("J22").Value = 731
$ws.Range("K22").Value = 3254.4443
$ws.Range("L22").Value = 731
$ws.Range("M22").Value = -2959.4443
$ws.Range("N22").Value = -1321

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2623.5833
$ws.Range("I27").Value = 3254.4443
$ws.Range("J27").Value = 731
$ws.Range("K27").Value = 3254.4443
$ws.Range("L27").Value = 731
$ws.Range("M27").Value = -3147.4443
$ws.Range("N27").Value = -945

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2376.2593
$ws.Range("I93").Value = 2335
$ws.Range("J93").Value = 2404.625
$ws.Range("K93").Value = 2335
$ws.Range("L93").Value = 2404.625
$ws.Range("M93").Value = -1087
$ws.Range("N93").Value = -4900.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3031.8333
$ws.Range("J100").Value = 3798.1
$ws.Range("L100").Value = 3798.1
$ws.Range("N100").Value = -4880.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5158.6665
$ws.Range("I136").Value = 1650.6666
$ws.Range("J136").Value = 8666.666999999999
$ws.Range("K136").Value = 4951.9998
$ws.Range("L136").Value = 26000.001
$ws.Range("M136").Value = -2401.9998
$ws.Range("N136").Value = -31100.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 84996.664
$ws.Range("J140").Value = 84996.664
$ws.Range("L140").Value = 84996.664
$ws.Range("N140").Value = -95356.664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28849222
$ws.Range("I132").Value = 40541584
$ws.Range("J132").Value = 8068.933
$ws.Range("K132").Value = 121624752
$ws.Range("L132").Value = 24206.799
$ws.Range("M132").Value = -121622222
$ws.Range("N132").Value = -29266.799
